$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("G2").Value = 88.72291666666666
$ws.Range("H2").Value = 266.16875
$ws.Range("I2").Value = 0.7675060578750151
$ws.Range("J2").Value = 0.7675060578750152
$ws.Range("K2").Value = 3.0
$ws.Range("M2").Value = 24.42119933333333
$ws.Range("N2").Value = 73.263598
$ws.Range("O2").Value = 0.4086816635579248
$ws.Range("P2").Value = 0.4086816635579248
$ws.Range("Q2").Value = 2166.720033351389
$ws.Range("R2").Value = 19500.4803001625
$ws.Range("S2").Value = 0.3136656525231462
$ws.Range("T2").Value = 0.3136656525231462

$ws.Range("E3").Value = 3.0
$ws.Range("G3").Value = 88.72291666666666
$ws.Range("H3").Value = 266.16875
$ws.Range("I3").Value = 0.7675060578750151
$ws.Range("J3").Value = 0.7675060578750152
$ws.Range("K3").Value = 3.0
$ws.Range("M3").Value = 33.48129
$ws.Range("N3").Value = 100.44387
$ws.Range("O3").Value = 0.5602996441124273
$ws.Range("P3").Value = 0.5602996441124273
$ws.Range("Q3").Value = 2970.5577025625
$ws.Range("R3").Value = 26735.0193230625
$ws.Range("S3").Value = 0.430033371081503
$ws.Range("T3").Value = 0.4300333710815031

$ws.Range("E4").Value = 3.0
$ws.Range("G4").Value = 88.72291666666666
$ws.Range("H4").Value = 266.16875
$ws.Range("I4").Value = 0.7675060578750151
$ws.Range("J4").Value = 0.7675060578750152
$ws.Range("K4").Value = 3.0
$ws.Range("M4").Value = 1.853554333333333
$ws.Range("N4").Value = 5.560663
$ws.Range("O4").Value = 0.03101869232964781
$ws.Range("P4").Value = 0.03101869232964781
$ws.Range("Q4").Value = 164.4527466534722
$ws.Range("R4").Value = 1480.07471988125
$ws.Range("S4").Value = 0.02380703427036596
$ws.Range("T4").Value = 0.02380703427036597

$ws.Range("E5").Value = 3.0
$ws.Range("G5").Value = 17.91585
$ws.Range("H5").Value = 53.74755
$ws.Range("I5").Value = 0.1549827702197958
$ws.Range("J5").Value = 0.1549827702197958
$ws.Range("K5").Value = 3.0
$ws.Range("M5").Value = 24.42119933333333
$ws.Range("N5").Value = 73.263598
$ws.Range("O5").Value = 0.4086816635579248
$ws.Range("P5").Value = 0.4086816635579248
$ws.Range("Q5").Value = 437.5265440761001
$ws.Range("R5").Value = 3937.7388966849
$ws.Range("S5").Value = 0.06333861635624177
$ws.Range("T5").Value = 0.06333861635624177

$ws.Range("E6").Value = 3.0
$ws.Range("G6").Value = 17.91585
$ws.Range("H6").Value = 53.74755
$ws.Range("I6").Value = 0.1549827702197958
$ws.Range("J6").Value = 0.1549827702197958
$ws.Range("K6").Value = 3.0
$ws.Range("M6").Value = 33.48129
$ws.Range("N6").Value = 100.44387
$ws.Range("O6").Value = 0.5602996441124273
$ws.Range("P6").Value = 0.5602996441124273
$ws.Range("Q6").Value = 599.8457694465001
$ws.Range("R6").Value = 5398.6119250185
$ws.Range("S6").Value = 0.08683679099770969
$ws.Range("T6").Value = 0.0868367909977097

$ws.Range("E7").Value = 3.0
$ws.Range("G7").Value = 17.91585
$ws.Range("H7").Value = 53.74755
$ws.Range("I7").Value = 0.1549827702197958
$ws.Range("J7").Value = 0.1549827702197958
$ws.Range("K7").Value = 3.0
$ws.Range("M7").Value = 1.853554333333333
$ws.Range("N7").Value = 5.560663
$ws.Range("O7").Value = 0.03101869232964781
$ws.Range("P7").Value = 0.03101869232964781
$ws.Range("Q7").Value = 33.20800140285
$ws.Range("R7").Value = 298.87201262565
$ws.Range("S7").Value = 0.00480736286584435
$ws.Range("T7").Value = 0.004807362865844351

$ws.Range("E8").Value = 3.0
$ws.Range("G8").Value = 8.960212333333333
$ws.Range("H8").Value = 26.880637
$ws.Range("I8").Value = 0.077511171905189
$ws.Range("J8").Value = 0.07751117190518901
$ws.Range("K8").Value = 3.0
$ws.Range("M8").Value = 24.42119933333333
$ws.Range("N8").Value = 73.263598
$ws.Range("O8").Value = 0.4086816635579248
$ws.Range("P8").Value = 0.4086816635579248
$ws.Range("Q8").Value = 218.8191314613251
$ws.Range("R8").Value = 1969.372183151926
$ws.Range("S8").Value = 0.03167739467853693
$ws.Range("T8").Value = 0.03167739467853693

$ws.Range("E9").Value = 3.0
$ws.Range("G9").Value = 8.960212333333333
$ws.Range("H9").Value = 26.880637
$ws.Range("I9").Value = 0.077511171905189
$ws.Range("J9").Value = 0.07751117190518901
$ws.Range("K9").Value = 3.0
$ws.Range("M9").Value = 33.48129
$ws.Range("N9").Value = 100.44387
$ws.Range("O9").Value = 0.5602996441124273
$ws.Range("P9").Value = 0.5602996441124273
$ws.Range("Q9").Value = 299.99946759391
$ws.Range("R9").Value = 2699.99520834519
$ws.Range("S9").Value = 0.04342948203321457
$ws.Range("T9").Value = 0.04342948203321458

$ws.Range("E10").Value = 3.0
$ws.Range("G10").Value = 8.960212333333333
$ws.Range("H10").Value = 26.880637
$ws.Range("I10").Value = 0.077511171905189
$ws.Range("J10").Value = 0.07751117190518901
$ws.Range("K10").Value = 3.0
$ws.Range("M10").Value = 1.853554333333333
$ws.Range("N10").Value = 5.560663
$ws.Range("O10").Value = 0.03101869232964781
$ws.Range("P10").Value = 0.03101869232964781
$ws.Range("Q10").Value = 16.60824039803678
$ws.Range("R10").Value = 149.474163582331
$ws.Range("S10").Value = 0.002404295193437499
$ws.Range("T10").Value = 0.002404295193437499
